# Edit script for quiz240930.xlsx
# Adds 31 new survey response rows (332-362) to the "Form_Responses1" table,
# expanding the table / used range from A1:N331 to A1:N362, and updates the
# sheet's selection to match the post-edit state (F372).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$rowData = @(
    ,(@(332, 45571.374858738425, 'p20236727@gmail.com', '인공지능융합학부', 20236727, '박진서', '민주 문자', '한글', '하나도 없다', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Red', '휴우, 그래도 반이나 남았네.', $null))
    ,(@(333, 45571.381704803236, 'vcx76613@gmail.com', '간호학과', 20246306, '황인태', '민주 문자', '한글', '하나도 없다', 0.9, '미국', '사회활동이나 자원활동에 덜 참여한다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(334, 45571.441933599534, 'tjtkdwns0909@naver.com', '체육학과', 20244123, '서상준', '민주 문자', '한글', '2개', 0.8, '대한민국', '건강이 좋지 않다', 'Red', '휴우, 그래도 반이나 남았네.', $null))
    ,(@(335, 45571.442870844912, 'chjames2005@naver.com', '데이터사이언스학부', 20243259, '최재현', '민주 문자', '한글', '2개', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Black', $null, '휴우, 그래도 반이나 남았네.'))
    ,(@(336, 45571.459191377318, 'doheehana@naver.com', '간호학과', 20246215, '김도희', '민주 문자', '한글', '1개', 0.8, '대한민국', '남들을 덜 신뢰한다', 'Black', $null, '휴우, 그래도 반이나 남았네.'))
    ,(@(337, 45571.475912708338, 'jyoon2233@naver.com', '소프트웨어학부', 20245190, '신정윤', '민주 문자', '한자', '하나도 없다', 0.9, '영국', '2배 정도 실직할 가능성이 높다', 'Red', '모름/기타', $null))
    ,(@(338, 45571.518326030091, 'dlaehdghks123@gmail.com', '금융재무학과', 20213533, '임동환', '엘리트 문자', '한자', '4개', 0.8, '대한민국', '2배 정도 실직할 가능성이 높다', 'Black', $null, '휴우, 그래도 반이나 남았네.'))
    ,(@(339, 45571.523208796294, '20242925@hallym.ac.kr', '경영학부', 20242925, '김민우', '민주 문자', '한글', '2개', 0.8, '대한민국', '남들을 덜 신뢰한다', 'Red', '휴우, 그래도 반이나 남았네.', $null))
    ,(@(340, 45571.53662082176, 'shasha4321@naver.com', '경영', 20213035, '정다영', '민주 문자', '한글', '하나도 없다', 0.2, '대한민국', '건강이 좋지 않다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(341, 45571.560315520837, 'sumine0601@naver.com', '광고홍보학과', 20202637, '장수민', '민주 문자', '한글', '하나도 없다', 0.8, '이탈리아', '건강이 좋지 않다', 'Black', $null, '휴우, 그래도 반이나 남았네.'))
    ,(@(342, 45571.560898402779, 'whrudghks030604@naver.com', '언어청각학부', 20243955, '조경환', '민주 문자', '한자', '하나도 없다', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(343, 45571.573620162038, 'sshee718@gmail.com', '환경생명공학과', 20243702, '권도운', '엘리트 문자', '한글', '2개', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(344, 45571.578612951387, 'lattace05@gmail.com', '정치행정학과', 20242439, '최동희', '민주 문자', '한글', '2개', 0.8, '대한민국', '남들을 덜 신뢰한다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(345, 45571.582787094907, 'r67890@naver.com', '스마트iot', 20205217, '이규형', '민주 문자', '한글', '하나도 없다', 0.5, '미국', '사회활동이나 자원활동에 덜 참여한다', 'Red', '휴우, 그래도 반이나 남았네.', $null))
    ,(@(346, 45571.596237905091, 'jiminn101777@gmail.com', '사회복지학부', 20242306, '권지민', '엘리트 문자', '한글', '1개', 0.5, '미국', '남들을 덜 신뢰한다', 'Black', $null, '휴우, 그래도 반이나 남았네.'))
    ,(@(347, 45571.600335300929, 'leedowon567@naver.com', '일본학과', 20211625, '이도원', '민주 문자', '한글', '하나도 없다', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Black', $null, '휴우, 그래도 반이나 남았네.'))
    ,(@(348, 45571.603621458329, 'dlakrp731@gmail.com', '콘텐츠it', 20195225, '이준수', '민주 문자', '한글', '하나도 없다', 0.9, '대한민국', '시간당 중위 임금이 60% 낮다', 'Red', '휴우, 그래도 반이나 남았네.', $null))
    ,(@(349, 45571.617213923615, 'simyenho8562@gmail.com', '경영대학', 20242981, '심연호', '민주 문자', '한글', '3개', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Red', '휴우, 그래도 반이나 남았네.', $null))
    ,(@(350, 45571.621937777774, 'twenty__dec@naver.com', '언어청각학부', 20243959, '채희주', '민주 문자', '한글', '하나도 없다', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(351, 45571.62493583333, 'sangim041113@gmail.com', '콘텐츠IT', 20235209, '용상임', '민주 문자', '한글', '하나도 없다', 0.8, '대한민국', '건강이 좋지 않다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(352, 45571.632109027778, 'wldus0859@gmail.com', '광고홍보학과', 20242618, '박지연', '민주 문자', '한글', '하나도 없다', 0.2, '대한민국', '남들을 덜 신뢰한다', 'Red', '휴우, 그래도 반이나 남았네.', $null))
    ,(@(353, 45571.645036527778, 'chetbaker22@naver.com', '철학전공', 20201031, '김채원', '민주 문자', '한글', '1개', 0.8, '대한민국', '건강이 좋지 않다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(354, 45571.648349930554, 'guj2205146@gmail.com', '심리학과', 20242101, '강의주', '민주 문자', '한글', '하나도 없다', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(355, 45571.655777685184, 'ystop061012@naver.com', '미래융합스쿨', 20246628, '손연수', '엘리트 문자', '한글', '하나도 없다', 0.9, '대한민국', '2배 정도 실직할 가능성이 높다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(356, 45571.663920127314, 'yuvin0612@naver.com', '생명과학과', 20203537, '장유빈', '민주 문자', '한글', '1개', 0.8, '대한민국', '남들을 덜 신뢰한다', 'Red', '휴우, 그래도 반이나 남았네.', $null))
    ,(@(357, 45571.667273530096, 'hyeonjin0976@gmail.com', '법학과', 20192717, '류현진', '민주 문자', '한글', '하나도 없다', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Black', $null, '휴우, 그래도 반이나 남았네.'))
    ,(@(358, 45571.693092476853, 'jangjunhyeok1001@naver.com', '반도체디스플레이스쿨', 20193341, '장준혁', '민주 문자', '한자', '1개', 0.2, '대한민국', '건강이 좋지 않다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(359, 45571.706114675922, 'jjy021026@gmail.com', '체육학과', 20214143, '전지환', '민주 문자', '한글', '1개', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Black', $null, '헐, 반 밖에 안 남았네.'))
    ,(@(360, 45571.715269675929, 'applehanul@naver.com', '화학과', 20233412, '신하늘', '민주 문자', '한글', '3개', 0.8, '대한민국', '시간당 중위 임금이 60% 낮다', 'Black', $null, '휴우, 그래도 반이나 남았네.'))
    ,(@(361, 45571.734542199076, 'ncu11069@naver.com', '이규민', 20231622, '이규민', '민주 문자', '한글', '하나도 없다', 0.8, '대한민국', '2배 정도 실직할 가능성이 높다', 'Red', '휴우, 그래도 반이나 남았네.', $null))
    ,(@(362, 45571.738939953706, 'buj5195193@naver.com', '법학과', 20192721, '백유진', '민주 문자', '한글', '2개', 0.8, '대한민국', '사회활동이나 자원활동에 덜 참여한다', 'Red', '휴우, 그래도 반이나 남았네.', $null))
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    # Grow the table by one row for every new response.
    $null = $lo.ListRows.Add()

    $ws.Cells.Item($r, 1).Value  = $entry[1]   # A: timestamp
    $ws.Cells.Item($r, 2).Value  = $entry[2]   # B: email address
    $ws.Cells.Item($r, 3).Value  = $entry[3]   # C: department
    $ws.Cells.Item($r, 4).Value  = $entry[4]   # D: student id
    $ws.Cells.Item($r, 5).Value  = $entry[5]   # E: name
    $ws.Cells.Item($r, 6).Value  = $entry[6]   # F: Q1
    $ws.Cells.Item($r, 7).Value  = $entry[7]   # G: Q2
    $ws.Cells.Item($r, 8).Value  = $entry[8]   # H: Q3
    $ws.Cells.Item($r, 9).Value  = $entry[9]   # I: Q4
    $ws.Cells.Item($r, 10).Value = $entry[10]  # J: Q5
    $ws.Cells.Item($r, 11).Value = $entry[11]  # K: Q6
    $ws.Cells.Item($r, 12).Value = $entry[12]  # L: Red/Black
    if ($entry[13] -ne $null) {
        $ws.Cells.Item($r, 13).Value = $entry[13]  # M: Q7 (branch 1)
    }
    if ($entry[14] -ne $null) {
        $ws.Cells.Item($r, 14).Value = $entry[14]  # N: Q7 (branch 2)
    }
}

# Reproduce the final saved selection / scroll state from the edit.
$null = $ws.Range("A321").Select()
$excel.ActiveWindow.ScrollRow = 321
$null = $ws.Range("F372").Select()
